# Updates Price (D) and Volume(1h) (E) columns for the cryptos sheet.
# Values are forced to text format ("@") to match the source data which
# stores these as literal strings (e.g. "501.40", "8.00") rather than numbers,
# preventing Excel from stripping meaningful trailing zeros or reformatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.606.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.735.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.734.48"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.81"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.354.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.732.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.632.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "501.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.12"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.97"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.89"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000136"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +8.90%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.09"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.94%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +13.45%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "439.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "49.78"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.950.77"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.05"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.80%  "
